# Auto-generated edit script applying numeric updates to Titan_Profits workbook
# Updates currentAveragePrice / Leve price / profit columns (H-N) on multiple sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(18, 8).Value = 841.3333  # H18: 697.9 -> 841.3333
$ws.Cells.Item(18, 9).Value = 623  # I18: 586.25 -> 623
$ws.Cells.Item(18, 10).Value = 950.5  # J18: 772.3333 -> 950.5
$ws.Cells.Item(18, 11).Value = 623  # K18: 586.25 -> 623
$ws.Cells.Item(18, 12).Value = 950.5  # L18: 772.3333 -> 950.5
$ws.Cells.Item(18, 13).Value = -339  # M18: -302.25 -> -339
$ws.Cells.Item(18, 14).Value = -1518.5  # N18: -1340.3333 -> -1518.5
$ws.Cells.Item(38, 8).Value = 362.8125  # H38: 372.5 -> 362.8125
$ws.Cells.Item(38, 10).Value = 640  # J38: 812.5 -> 640
$ws.Cells.Item(38, 12).Value = 1920  # L38: 2437.5 -> 1920
$ws.Cells.Item(38, 14).Value = -2664  # N38: -3181.5 -> -2664
$ws.Cells.Item(40, 8).Value = 2286.1667  # H40: 2233.5386 -> 2286.1667
$ws.Cells.Item(40, 10).Value = 2699.6667  # J40: 2425.25 -> 2699.6667
$ws.Cells.Item(40, 12).Value = 2699.6667  # L40: 2425.25 -> 2699.6667
$ws.Cells.Item(40, 14).Value = -3049.6667  # N40: -2775.25 -> -3049.6667
$ws.Cells.Item(43, 8).Value = 840.2857  # H43: 880.3333 -> 840.2857
$ws.Cells.Item(43, 9).Value = 800  # I43: 866.6667 -> 800
$ws.Cells.Item(43, 11).Value = 800  # K43: 866.6667 -> 800
$ws.Cells.Item(43, 13).Value = -731  # M43: -797.6667 -> -731
$ws.Cells.Item(96, 8).Value = 775  # H96: 627.2727 -> 775
$ws.Cells.Item(96, 9).Value = 640  # I96: 700 -> 640
$ws.Cells.Item(96, 10).Value = 1000  # J96: 611.1111 -> 1000
$ws.Cells.Item(96, 11).Value = 1920  # K96: 2100 -> 1920
$ws.Cells.Item(96, 12).Value = 3000  # L96: 1833.3333 -> 3000
$ws.Cells.Item(96, 13).Value = -547  # M96: -727 -> -547
$ws.Cells.Item(96, 14).Value = -5746  # N96: -4579.3333 -> -5746
$ws.Cells.Item(112, 8).Value = 9616396  # H112: 9260267 -> 9616396
$ws.Cells.Item(112, 10).Value = 9616396  # J112: 9260267 -> 9616396
$ws.Cells.Item(112, 12).Value = 28849188  # L112: 27780801 -> 28849188
$ws.Cells.Item(112, 14).Value = -28851404  # N112: -27783017 -> -28851404

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 1013.58826  # H2: 57081.11 -> 1013.58826
$ws.Cells.Item(2, 9).Value = 765.5714  # I2: 68044 -> 765.5714
$ws.Cells.Item(2, 10).Value = 2171  # J2: 2266.6667 -> 2171
$ws.Cells.Item(2, 11).Value = 765.5714  # K2: 68044 -> 765.5714
$ws.Cells.Item(2, 12).Value = 2171  # L2: 2266.6667 -> 2171
$ws.Cells.Item(2, 13).Value = -652.5714  # M2: -67931 -> -652.5714
$ws.Cells.Item(2, 14).Value = -2397  # N2: -2492.6667 -> -2397
$ws.Cells.Item(45, 8).Value = 1099  # H45: 1115.9375 -> 1099
$ws.Cells.Item(45, 9).Value = 1036.5  # I45: 1060.4166 -> 1036.5
$ws.Cells.Item(45, 10).Value = 1224  # J45: 1282.5 -> 1224
$ws.Cells.Item(45, 11).Value = 1036.5  # K45: 1060.4166 -> 1036.5
$ws.Cells.Item(45, 12).Value = 1224  # L45: 1282.5 -> 1224
$ws.Cells.Item(45, 13).Value = -659.5  # M45: -683.4166 -> -659.5
$ws.Cells.Item(45, 14).Value = -1978  # N45: -2036.5 -> -1978
$ws.Cells.Item(61, 8).Value = 2314.2  # H61: 2861.8262 -> 2314.2
$ws.Cells.Item(61, 9).Value = 1656.7273  # I61: 2016.2667 -> 1656.7273
$ws.Cells.Item(61, 10).Value = 4122.25  # J61: 4447.25 -> 4122.25
$ws.Cells.Item(61, 11).Value = 1656.7273  # K61: 2016.2667 -> 1656.7273
$ws.Cells.Item(61, 12).Value = 4122.25  # L61: 4447.25 -> 4122.25
$ws.Cells.Item(61, 13).Value = -1444.7273  # M61: -1804.2667 -> -1444.7273
$ws.Cells.Item(61, 14).Value = -4546.25  # N61: -4871.25 -> -4546.25
$ws.Cells.Item(116, 8).Value = 1013.58826  # H116: 57081.11 -> 1013.58826
$ws.Cells.Item(116, 9).Value = 765.5714  # I116: 68044 -> 765.5714
$ws.Cells.Item(116, 10).Value = 2171  # J116: 2266.6667 -> 2171
$ws.Cells.Item(116, 11).Value = 765.5714  # K116: 68044 -> 765.5714
$ws.Cells.Item(116, 12).Value = 2171  # L116: 2266.6667 -> 2171
$ws.Cells.Item(116, 13).Value = 1528.4286  # M116: -65750 -> 1528.4286
$ws.Cells.Item(116, 14).Value = -6759  # N116: -6854.6667 -> -6759
$ws.Cells.Item(136, 8).Value = 2314.2  # H136: 2861.8262 -> 2314.2
$ws.Cells.Item(136, 9).Value = 1656.7273  # I136: 2016.2667 -> 1656.7273
$ws.Cells.Item(136, 10).Value = 4122.25  # J136: 4447.25 -> 4122.25
$ws.Cells.Item(136, 11).Value = 4970.1819  # K136: 6048.800099999999 -> 4970.1819
$ws.Cells.Item(136, 12).Value = 12366.75  # L136: 13341.75 -> 12366.75
$ws.Cells.Item(136, 13).Value = -2420.1819  # M136: -3498.800099999999 -> -2420.1819
$ws.Cells.Item(136, 14).Value = -17466.75  # N136: -18441.75 -> -17466.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 1013.58826  # H3: 57081.11 -> 1013.58826
$ws.Cells.Item(3, 9).Value = 765.5714  # I3: 68044 -> 765.5714
$ws.Cells.Item(3, 10).Value = 2171  # J3: 2266.6667 -> 2171
$ws.Cells.Item(3, 11).Value = 765.5714  # K3: 68044 -> 765.5714
$ws.Cells.Item(3, 12).Value = 2171  # L3: 2266.6667 -> 2171
$ws.Cells.Item(3, 13).Value = -651.5714  # M3: -67930 -> -651.5714
$ws.Cells.Item(3, 14).Value = -2399  # N3: -2494.6667 -> -2399

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 1418.7021  # H31: 1455.2291 -> 1418.7021
$ws.Cells.Item(31, 9).Value = 943.7  # I31: 933.43335 -> 943.7
$ws.Cells.Item(31, 10).Value = 2256.9412  # J31: 2324.889 -> 2256.9412
$ws.Cells.Item(31, 11).Value = 943.7  # K31: 933.43335 -> 943.7
$ws.Cells.Item(31, 12).Value = 2256.9412  # L31: 2324.889 -> 2256.9412
$ws.Cells.Item(31, 13).Value = -648.7  # M31: -638.43335 -> -648.7
$ws.Cells.Item(31, 14).Value = -2846.9412  # N31: -2914.889 -> -2846.9412
$ws.Cells.Item(33, 8).Value = 26125  # H33: 33010.332 -> 26125
$ws.Cells.Item(33, 9).Value = 26125  # I33: 33010.332 -> 26125
$ws.Cells.Item(33, 11).Value = 26125  # K33: 33010.332 -> 26125
$ws.Cells.Item(33, 13).Value = -25746  # M33: -32631.332 -> -25746
$ws.Cells.Item(34, 8).Value = 1418.7021  # H34: 1455.2291 -> 1418.7021
$ws.Cells.Item(34, 9).Value = 943.7  # I34: 933.43335 -> 943.7
$ws.Cells.Item(34, 10).Value = 2256.9412  # J34: 2324.889 -> 2256.9412
$ws.Cells.Item(34, 11).Value = 943.7  # K34: 933.43335 -> 943.7
$ws.Cells.Item(34, 12).Value = 2256.9412  # L34: 2324.889 -> 2256.9412
$ws.Cells.Item(34, 13).Value = -741.7  # M34: -731.43335 -> -741.7
$ws.Cells.Item(34, 14).Value = -2660.9412  # N34: -2728.889 -> -2660.9412
$ws.Cells.Item(125, 8).Value = 0  # H125: 54900 -> 0
$ws.Cells.Item(125, 10).Value = 0  # J125: 54900 -> 0
$ws.Cells.Item(125, 12).Value = 0  # L125: 54900 -> 0
$ws.Cells.Item(125, 14).ClearContents()  # N125: remove cell (was -59820)
$ws.Cells.Item(134, 8).Value = 2595.9688  # H134: 2930.3447 -> 2595.9688
$ws.Cells.Item(134, 9).Value = 1377.96  # I134: 1644.3684 -> 1377.96
$ws.Cells.Item(134, 10).Value = 6946  # J134: 5373.7 -> 6946
$ws.Cells.Item(134, 11).Value = 4133.88  # K134: 4933.1052 -> 4133.88
$ws.Cells.Item(134, 12).Value = 20838  # L134: 16121.1 -> 20838
$ws.Cells.Item(134, 13).Value = -1598.88  # M134: -2398.1052 -> -1598.88
$ws.Cells.Item(134, 14).Value = -25908  # N134: -21191.1 -> -25908

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(117, 8).Value = 998.75  # H117: 1140 -> 998.75
$ws.Cells.Item(117, 9).Value = 1031.6  # I117: 1233.3334 -> 1031.6
$ws.Cells.Item(117, 10).Value = 944  # J117: 1000 -> 944
$ws.Cells.Item(117, 11).Value = 3094.8  # K117: 3700.0002 -> 3094.8
$ws.Cells.Item(117, 12).Value = 2832  # L117: 3000 -> 2832
$ws.Cells.Item(117, 13).Value = 347.2000000000003  # M117: -258.0001999999999 -> 347.2000000000003
$ws.Cells.Item(117, 14).Value = -9716  # N117: -9884 -> -9716
$ws.Cells.Item(122, 8).Value = 1076.1428  # H122: 1109 -> 1076.1428
$ws.Cells.Item(122, 10).Value = 1076.1428  # J122: 1109 -> 1076.1428
$ws.Cells.Item(122, 12).Value = 9685.2852  # L122: 9981 -> 9685.2852
$ws.Cells.Item(122, 14).Value = -14585.2852  # N122: -14881 -> -14585.2852
$ws.Cells.Item(131, 8).Value = 1449.037  # H131: 1399.322 -> 1449.037
$ws.Cells.Item(131, 9).Value = 342.66666  # I131: 338.1 -> 342.66666
$ws.Cells.Item(131, 10).Value = 1670.3112  # J131: 1615.898 -> 1670.3112
$ws.Cells.Item(131, 11).Value = 1027.99998  # K131: 1014.3 -> 1027.99998
$ws.Cells.Item(131, 12).Value = 5010.9336  # L131: 4847.694 -> 5010.9336
$ws.Cells.Item(131, 13).Value = 4012.00002  # M131: 4025.7 -> 4012.00002
$ws.Cells.Item(131, 14).Value = -15090.9336  # N131: -14927.694 -> -15090.9336
$ws.Cells.Item(134, 8).Value = 2700.647  # H134: 3889.6843 -> 2700.647
$ws.Cells.Item(134, 9).Value = 1798.5333  # I134: 2344.7334 -> 1798.5333
$ws.Cells.Item(134, 10).Value = 9466.5  # J134: 9683.25 -> 9466.5
$ws.Cells.Item(134, 11).Value = 5395.5999  # K134: 7034.2002 -> 5395.5999
$ws.Cells.Item(134, 12).Value = 28399.5  # L134: 29049.75 -> 28399.5
$ws.Cells.Item(134, 13).Value = -325.5999000000002  # M134: -1964.2002 -> -325.5999000000002
$ws.Cells.Item(134, 14).Value = -38539.5  # N134: -39189.75 -> -38539.5
$ws.Cells.Item(140, 8).Value = 4569.9116  # H140: 4577.5293 -> 4569.9116
$ws.Cells.Item(140, 9).Value = 6099.579  # I140: 6113.2104 -> 6099.579
$ws.Cells.Item(140, 11).Value = 18298.737  # K140: 18339.6312 -> 18298.737
$ws.Cells.Item(140, 13).Value = -13118.737  # M140: -13159.6312 -> -13118.737
$ws.Cells.Item(141, 8).Value = 4523.9  # H141: 4570.9 -> 4523.9
$ws.Cells.Item(141, 9).Value = 5391.2856  # I141: 6034.8335 -> 5391.2856
$ws.Cells.Item(141, 10).Value = 2500  # J141: 2375 -> 2500
$ws.Cells.Item(141, 11).Value = 16173.8568  # K141: 18104.5005 -> 16173.8568
$ws.Cells.Item(141, 12).Value = 7500  # L141: 7125 -> 7500
$ws.Cells.Item(141, 13).Value = -10993.8568  # M141: -12924.5005 -> -10993.8568
$ws.Cells.Item(141, 14).Value = -17860  # N141: -17485 -> -17860

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(113, 8).Value = 1257.24  # H113: 1548.5652 -> 1257.24
$ws.Cells.Item(113, 9).Value = 1151.6471  # I113: 1519.8235 -> 1151.6471
$ws.Cells.Item(113, 10).Value = 1481.625  # J113: 1630 -> 1481.625
$ws.Cells.Item(113, 11).Value = 1151.6471  # K113: 1519.8235 -> 1151.6471
$ws.Cells.Item(113, 12).Value = 1481.625  # L113: 1630 -> 1481.625
$ws.Cells.Item(113, 13).Value = 1018.3529  # M113: 650.1765 -> 1018.3529
$ws.Cells.Item(113, 14).Value = -5821.625  # N113: -5970 -> -5821.625
$ws.Cells.Item(122, 8).Value = 5556805.5  # H122: 2778980.5 -> 5556805.5
$ws.Cells.Item(122, 9).Value = 11111111  # I122: 3704474 -> 11111111
$ws.Cells.Item(122, 11).Value = 33333333  # K122: 11113422 -> 33333333
$ws.Cells.Item(122, 13).Value = -33330883  # M122: -11110972 -> -33330883
$ws.Cells.Item(132, 8).Value = 4351.1177  # H132: 4433.2646 -> 4351.1177
$ws.Cells.Item(132, 9).Value = 4553.1904  # I132: 4519.909 -> 4553.1904
$ws.Cells.Item(132, 10).Value = 4024.6924  # J132: 4274.4165 -> 4024.6924
$ws.Cells.Item(132, 11).Value = 13659.5712  # K132: 13559.727 -> 13659.5712
$ws.Cells.Item(132, 12).Value = 12074.0772  # L132: 12823.2495 -> 12074.0772
$ws.Cells.Item(132, 13).Value = -11129.5712  # M132: -11029.727 -> -11129.5712
$ws.Cells.Item(132, 14).Value = -17134.0772  # N132: -17883.2495 -> -17134.0772

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 842.75  # H22: 1258.3334 -> 842.75
$ws.Cells.Item(22, 9).Value = 218.4  # I22: 450 -> 218.4
$ws.Cells.Item(22, 10).Value = 1883.3334  # J22: 1662.5 -> 1883.3334
$ws.Cells.Item(22, 11).Value = 218.4  # K22: 450 -> 218.4
$ws.Cells.Item(22, 12).Value = 1883.3334  # L22: 1662.5 -> 1883.3334
$ws.Cells.Item(22, 13).Value = 76.59999999999999  # M22: -155 -> 76.59999999999999
$ws.Cells.Item(22, 14).Value = -2473.3334  # N22: -2252.5 -> -2473.3334
$ws.Cells.Item(27, 8).Value = 842.75  # H27: 1258.3334 -> 842.75
$ws.Cells.Item(27, 9).Value = 218.4  # I27: 450 -> 218.4
$ws.Cells.Item(27, 10).Value = 1883.3334  # J27: 1662.5 -> 1883.3334
$ws.Cells.Item(27, 11).Value = 218.4  # K27: 450 -> 218.4
$ws.Cells.Item(27, 12).Value = 1883.3334  # L27: 1662.5 -> 1883.3334
$ws.Cells.Item(27, 13).Value = -111.4  # M27: -343 -> -111.4
$ws.Cells.Item(27, 14).Value = -2097.3334  # N27: -1876.5 -> -2097.3334
$ws.Cells.Item(40, 8).Value = 3127.9443  # H40: 3109.3157 -> 3127.9443
$ws.Cells.Item(40, 9).Value = 1825.75  # I40: 1899.6666 -> 1825.75
$ws.Cells.Item(40, 10).Value = 3500  # J40: 3336.125 -> 3500
$ws.Cells.Item(40, 11).Value = 1825.75  # K40: 1899.6666 -> 1825.75
$ws.Cells.Item(40, 12).Value = 3500  # L40: 3336.125 -> 3500
$ws.Cells.Item(40, 13).Value = -1689.75  # M40: -1763.6666 -> -1689.75
$ws.Cells.Item(40, 14).Value = -3772  # N40: -3608.125 -> -3772
$ws.Cells.Item(136, 8).Value = 4726.069  # H136: 4187.1763 -> 4726.069
$ws.Cells.Item(136, 9).Value = 2214.9524  # I136: 1993.1538 -> 2214.9524
$ws.Cells.Item(136, 11).Value = 6644.8572  # K136: 5979.4614 -> 6644.8572
$ws.Cells.Item(136, 13).Value = -4094.8572  # M136: -3429.4614 -> -4094.8572

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 14305484  # H62: 16689069 -> 14305484
$ws.Cells.Item(62, 9).Value = 20023256  # I62: 25028076 -> 20023256
$ws.Cells.Item(62, 11).Value = 20023256  # K62: 25028076 -> 20023256
$ws.Cells.Item(62, 13).Value = -20022632  # M62: -25027452 -> -20022632
$ws.Cells.Item(65, 8).Value = 14305484  # H65: 16689069 -> 14305484
$ws.Cells.Item(65, 9).Value = 20023256  # I65: 25028076 -> 20023256
$ws.Cells.Item(65, 11).Value = 100116280  # K65: 125140380 -> 100116280
$ws.Cells.Item(65, 13).Value = -100113160  # M65: -125137260 -> -100113160
$ws.Cells.Item(122, 8).Value = 23654.912  # H122: 28472.475 -> 23654.912
$ws.Cells.Item(122, 9).Value = 30847.295  # I122: 34927.367 -> 30847.295
$ws.Cells.Item(122, 10).Value = 3276.5  # J122: 4266.625 -> 3276.5
$ws.Cells.Item(122, 11).Value = 92541.88499999999  # K122: 104782.101 -> 92541.88499999999
$ws.Cells.Item(122, 12).Value = 9829.5  # L122: 12799.875 -> 9829.5
$ws.Cells.Item(122, 13).Value = -90091.88499999999  # M122: -102332.101 -> -90091.88499999999
$ws.Cells.Item(122, 14).Value = -14729.5  # N122: -17699.875 -> -14729.5
